# Auto-generated from the OOXML diff.
# All changed cells (D2:E51 price/volume table) are plain inline strings in the
# source workbook. Excel's COM layer auto-coerces numeric-looking text (e.g.
# "22.00", "1.20", "0.0000119") into numbers when assigned via .Value, which would
# silently drop significant trailing/leading zeros and change the cell type from
# text to number. Forcing NumberFormat = "@" (Text) before the assignment keeps
# every written value exactly as authored in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.251.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.177.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -8.29%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.67"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.01"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.96%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.607"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.177.03"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -8.37%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.95%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.16%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.723.83"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -8.42%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -9.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.250.05"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.172.92"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -8.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.74"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.98%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "353.14"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.18"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.67%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.62"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -6.50%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000119"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.54%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.504"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.59"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.74%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.53"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -7.41%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -7.34%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.20"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.58%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.63"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.19%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "154.06"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.27%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.05"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.54%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.49"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.615.27"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.18"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.37"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.98"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -7.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0650"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.44%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "319.70"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -6.68%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -7.83%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.07%  "
